$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.384.26'
$ws.Range("E2").Value = '  -0.77%  '
$ws.Range("D3").Value = '1.638.82'
$ws.Range("E3").Value = '  -1.65%  '
$ws.Range("E4").Value = '  +0.03%  '
$ws.Range("D5").Value = '''211.65'
$ws.Range("E5").Value = '  -1.54%  '
$ws.Range("D6").Value = '''0.528'
$ws.Range("E6").Value = '  +3.89%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '''22.93'
$ws.Range("E8").Value = '  -3.68%  '
$ws.Range("E9").Value = '  -1.99%  '
$ws.Range("D11").Value = '''0.0893'
$ws.Range("E11").Value = '  +1.67%  '
$ws.Range("D12").Value = '1.871.11'
$ws.Range("E12").Value = '  -1.63%  '
$ws.Range("D13").Value = '1.636.93'
$ws.Range("E13").Value = '  -1.70%  '
$ws.Range("E14").Value = '  -2.51%  '
$ws.Range("E15").Value = '  +0.20%  '
$ws.Range("D16").Value = '''64.33'
$ws.Range("E16").Value = '  -2.97%  '
$ws.Range("D17").Value = '27.357.84'
$ws.Range("E17").Value = '  -0.80%  '
$ws.Range("D18").Value = '''228.82'
$ws.Range("E18").Value = '  -5.86%  '
$ws.Range("E19").Value = '  -1.28%  '
$ws.Range("D20").Value = '''7.55'
$ws.Range("E20").Value = '  -1.09%  '
$ws.Range("E21").Value = '  -0.01%  '
$ws.Range("D23").Value = '''9.59'
$ws.Range("E23").Value = '  +3.16%  '
$ws.Range("E24").Value = '  -0.39%  '
$ws.Range("D25").Value = '''146.97'
$ws.Range("E25").Value = '  -0.01%  '
$ws.Range("E26").Value = '  -3.22%  '
$ws.Range("D27").Value = '''0.112'
$ws.Range("E27").Value = '  +1.08%  '
$ws.Range("E28").Value = '  -0.02%  '
$ws.Range("D29").Value = '''15.51'
$ws.Range("E29").Value = '  -5.71%  '
$ws.Range("E30").Value = '  -3.94%  '
$ws.Range("D31").Value = '''0.0483'
$ws.Range("E31").Value = '  -3.61%  '
$ws.Range("D32").Value = '''3.27'
$ws.Range("E32").Value = '  -2.23%  '
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("D34").Value = '1.412.15'
$ws.Range("E34").Value = '  -3.85%  '
$ws.Range("E35").Value = '  +0.51%  '
$ws.Range("E36").Value = '  -0.30%  '
$ws.Range("E37").Value = '  -1.69%  '
$ws.Range("D38").Value = '''0.879'
$ws.Range("E38").Value = '  -5.32%  '
$ws.Range("E39").Value = '  -3.10%  '
$ws.Range("E40").Value = '  +1.12%  '
$ws.Range("E41").Value = '  +0.03%  '
$ws.Range("D42").Value = '''2.46'
$ws.Range("E42").Value = '  -1.95%  '
$ws.Range("E43").Value = '  +1.88%  '
$ws.Range("D44").Value = '''2.24'
$ws.Range("E44").Value = '  +0.58%  '
$ws.Range("D45").Value = '''0.791'
$ws.Range("E45").Value = '  +0.69%  '
$ws.Range("D46").Value = '''64.46'
$ws.Range("E46").Value = '  -7.15%  '
$ws.Range("D47").Value = '1.780.67'
$ws.Range("E47").Value = '  -1.57%  '
$ws.Range("D48").Value = '''1.66'
$ws.Range("E48").Value = '  -3.65%  '
$ws.Range("D49").Value = '''87.69'
$ws.Range("E49").Value = '  -1.84%  '
$ws.Range("E50").Value = '  -1.65%  '
